$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the City value in J2 (previously "Mehboobnagar")
$ws.Range("J2").ClearContents()

# Select J2 to match the resulting workbook view selection
$ws.Range("J2").Select()
